# Actualización automática del inventario, Google Sheets y productos.json
# Adds a new inventory row (row 35) for "Correa de plotter para impresora HP DesignJet"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35

$ws.Cells.Item($row, 1).Value = "O3F6AE"
$ws.Cells.Item($row, 2).Value = "Correa de plotter para impresora HP DesignJet"
$ws.Cells.Item($row, 3).Value = "500 510 800 815 de 42`""
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 350000
$ws.Cells.Item($row, 6).Value = 2
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E35-D35)*G35"
$ws.Cells.Item($row, 9).Formula = "=D35*F35"
$ws.Cells.Item($row, 10).Value = 0
